$wb = $excel.ActiveWorkbook

# --- Summary sheet: update name, income, totals, net worth, ratio ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = "Kamal Al Shehhi"
$wsSummary.Range("B4").Value = 2480.82
$wsSummary.Range("B6").Value = 4833
$wsSummary.Range("B7").Value = 9516
$wsSummary.Range("B8").Value = -4683
$wsSummary.Range("B9").Value = 0.51

# --- Assets sheet: remove the "Vehicles / Economy Car" row, leaving
#     "Liquid Assets / Savings Account" as the sole asset line, and update
#     its value along with the TOTAL ASSETS figure ---
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Rows("2").Delete()
$wsAssets.Range("C2").Value = 4833
$wsAssets.Range("C3").Value = 4833

# --- Liabilities sheet: remove the "Auto Loans / Vehicle Loan 1" row,
#     leaving "Credit Cards / Credit Card Balance" as the sole liability
#     line, and update its figures along with TOTAL LIABILITIES ---
$wsLiabilities = $wb.Worksheets.Item("Liabilities")
$wsLiabilities.Rows("2").Delete()
$wsLiabilities.Range("C2").Value = 9516
$wsLiabilities.Range("D2").Value = 476
$wsLiabilities.Range("E2").Value = 1
$wsLiabilities.Range("C3").Value = 9516
